# "changes in design all pages and Near By Hospital"
# Adds a new "HEALTH RECORDS" row (row 6) to Sheet1, mirroring the
# formatting of the other URL rows (center aligned, wrapped params column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Populate the new row's cells. Write the URL (C6) before the label (B6) so
# the shared-string table grows in URL, label, params order.
$ws.Range("C6").Value = "http://192.168.100.19/thaimaiapp/api/mother/mHealthRecord"
$ws.Range("B6").Value = "HEALTH RECORDS"
$ws.Range("D6").Value = "picmeId=1000000000001" + [char]10 + "mid=1"

# Match formatting used by the other data rows (B/C centered; D centered + wrapped).
$ws.Range("B6:C6").HorizontalAlignment = -4108
$ws.Range("B6:C6").VerticalAlignment = -4108

$ws.Range("D6").HorizontalAlignment = -4108
$ws.Range("D6").VerticalAlignment = -4108
$ws.Range("D6").WrapText = $true

$ws.Rows.Item(6).RowHeight = 30

# Leave the selection parked on the newly added URL cell.
[void]$ws.Range("C6").Select()
